$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1624649859943978
$ws.Range("C2").Value = 0.6274509803921569
$ws.Range("J2").Value = 0.02521008403361345
$ws.Range("P2").Value = 0.1204481792717087
$ws.Range("S2").Value = 0.06442577030812324
$ws.Range("B3").Value = 0.004405286343612335
$ws.Range("C3").Value = 0.00881057268722467
$ws.Range("J3").Value = 0.01762114537444934
$ws.Range("P3").Value = 0.775330396475771
$ws.Range("S3").Value = 0.1938325991189427
$ws.Range("J4").Value = 0.06818181818181818
$ws.Range("P4").Value = 0.5454545454545454
$ws.Range("S4").Value = 0.3863636363636364
$ws.Range("B6").Value = 0.08108108108108109
$ws.Range("F6").Value = 0.03783783783783784
$ws.Range("J6").Value = 0.3027027027027027
$ws.Range("O6").Value = 0.01621621621621622
$ws.Range("Q6").Value = 0.1783783783783784
$ws.Range("R6").Value = 0.06486486486486487
$ws.Range("S6").Value = 0.3189189189189189
$ws.Range("B7").Value = 0.1153846153846154
$ws.Range("D7").Value = 0.01098901098901099
$ws.Range("F7").Value = 0.04395604395604396
$ws.Range("J7").Value = 0.1373626373626374
$ws.Range("O7").Value = 0.005494505494505495
$ws.Range("Q7").Value = 0.2307692307692308
$ws.Range("R7").Value = 0.05494505494505494
$ws.Range("S7").Value = 0.4010989010989011
$ws.Range("B8").Value = 0.1236673773987207
$ws.Range("D8").Value = 0.01066098081023454
$ws.Range("E8").Value = 0.002132196162046908
$ws.Range("F8").Value = 0.0511727078891258
$ws.Range("J8").Value = 0.1364605543710021
$ws.Range("O8").Value = 0.0255863539445629
$ws.Range("Q8").Value = 0.1812366737739872
$ws.Range("R8").Value = 0.07462686567164178
$ws.Range("S8").Value = 0.394456289978678
$ws.Range("B9").Value = 0.1155778894472362
$ws.Range("D9").Value = 0.02512562814070352
$ws.Range("F9").Value = 0.03015075376884422
$ws.Range("J9").Value = 0.1306532663316583
$ws.Range("O9").Value = 0.01507537688442211
$ws.Range("Q9").Value = 0.2010050251256282
$ws.Range("R9").Value = 0.06532663316582915
$ws.Range("S9").Value = 0.4170854271356784
$ws.Range("B10").Value = 0.1316779533483823
$ws.Range("D10").Value = 0.02407825432656132
$ws.Range("F10").Value = 0.0654627539503386
$ws.Range("J10").Value = 0.145974416854778
$ws.Range("O10").Value = 0.01354401805869074
$ws.Range("Q10").Value = 0.2302483069977427
$ws.Range("R10").Value = 0.06696764484574869
$ws.Range("S10").Value = 0.3220466516177577
$ws.Range("G11").Value = 0.1304347826086956
$ws.Range("J11").Value = 0.09420289855072464
$ws.Range("K11").Value = 0.2028985507246377
$ws.Range("L11").Value = 0.5543478260869565
$ws.Range("S11").Value = 0.01811594202898551
$ws.Range("G12").Value = 0.7453416149068323
$ws.Range("J12").Value = 0.1863354037267081
$ws.Range("K12").Value = 0.006211180124223602
$ws.Range("L12").Value = 0.04347826086956522
$ws.Range("S12").Value = 0.01863354037267081
$ws.Range("G13").Value = 0.7111111111111111
$ws.Range("J13").Value = 0.2444444444444444
$ws.Range("S13").Value = 0.04444444444444445
$ws.Range("F15").Value = 0.01744186046511628
$ws.Range("H15").Value = 0.1453488372093023
$ws.Range("I15").Value = 0.1104651162790698
$ws.Range("J15").Value = 0.3488372093023256
$ws.Range("K15").Value = 0.04651162790697674
$ws.Range("M15").Value = 0.005813953488372093
$ws.Range("O15").Value = 0.04069767441860465
$ws.Range("S15").Value = 0.2848837209302326
$ws.Range("F16").Value = 0.0131578947368421
$ws.Range("H16").Value = 0.2105263157894737
$ws.Range("I16").Value = 0.1096491228070175
$ws.Range("J16").Value = 0.3991228070175439
$ws.Range("K16").Value = 0.07456140350877193
$ws.Range("M16").Value = 0.01754385964912281
$ws.Range("O16").Value = 0.05263157894736842
$ws.Range("S16").Value = 0.1228070175438596
$ws.Range("F17").Value = 0.01397205588822355
$ws.Range("H17").Value = 0.2035928143712575
$ws.Range("I17").Value = 0.08782435129740519
$ws.Range("J17").Value = 0.4231536926147705
$ws.Range("K17").Value = 0.07584830339321358
$ws.Range("M17").Value = 0.01796407185628742
$ws.Range("O17").Value = 0.03592814371257485
$ws.Range("S17").Value = 0.1417165668662675
$ws.Range("H18").Value = 0.1962025316455696
$ws.Range("I18").Value = 0.1075949367088608
$ws.Range("J18").Value = 0.4240506329113924
$ws.Range("K18").Value = 0.120253164556962
$ws.Range("M18").Value = 0.006329113924050633
$ws.Range("O18").Value = 0.03164556962025317
$ws.Range("S18").Value = 0.1139240506329114
$ws.Range("F19").Value = 0.009819967266775777
$ws.Range("H19").Value = 0.220949263502455
$ws.Range("I19").Value = 0.07774140752864157
$ws.Range("J19").Value = 0.3772504091653028
$ws.Range("K19").Value = 0.1104746317512275
$ws.Range("M19").Value = 0.02700490998363339
$ws.Range("N19").Value = 0.001636661211129296
$ws.Range("O19").Value = 0.05646481178396072
$ws.Range("S19").Value = 0.118657937806874